# VERSION 0.9 ORM Edition
#
# Rewrites the "Реализация:" paragraph with the new bot description,
# moves the _GoBack bookmark into that paragraph, and reflows the
# grammar-check markup in the "Используемые технологии:" paragraph
# (splits the "фантазия(" run out of the "фантазия(ее " run and wraps
# it in a gramStart/gramEnd pair). The trailing bookmark that used to
# sit at the end of the "Скриншоты" paragraph is removed since it now
# lives earlier in the document.

$d = $word.ActiveDocument

# Courier New run formatting shared by every run in this document body.
$rPr = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr>'

function Get-ParagraphRange([string]$needle) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph containing '$needle' was not found"
    }
    $rng.Expand(4) | Out-Null   # wdParagraph -> grow to the whole paragraph, pilcrow included
    return $rng
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. "Реализация: ..." paragraph -----------------------------------
$implementationXml = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val=`"a3`"/>$rPr</w:pPr>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Реализация: </w:t></w:r>" +
    "<w:r>$rPr<w:t>мой проект это бот имеющий</w:t></w:r>" +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> несколько различных функций, главой особенностью является свободная база знаний которую могут пополнять все желающие, так же вы можете узнать погоду в своем городе, бросить кубик, или засечь время таймером </w:t></w:r>" +
    "</w:p>"

$implementationRange = Get-ParagraphRange("Реализация:")
$implementationRange.InsertXML($implementationXml)

# --- 2. "Используемые технологии: ..." paragraph -----------------------
$technologiesXml = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val=`"a3`"/>$rPr</w:pPr>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Используемые технологии: </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r>$rPr<w:t>Python</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">, </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r>$rPr<w:t>telegram.ext</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">, </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r>$rPr<w:t>кривыее</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> руки, </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/><w:r>$rPr<w:t>фантазия(</w:t></w:r><w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">ее </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r>$rPr<w:t>отсуцтвие</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r>$rPr<w:t>) и активное использование приемов ООП и ЧЖ</w:t></w:r>" +
    "</w:p>"

$technologiesRange = Get-ParagraphRange("Используемые технологии:")
$technologiesRange.InsertXML($technologiesXml)

# --- 3. "Скриншоты" paragraph: drop the now-relocated bookmark ---------
$screenshotsXml = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val=`"a3`"/>$rPr</w:pPr>" +
    "<w:r>$rPr<w:t>Ск</w:t></w:r>" +
    "<w:r>$rPr<w:t>риншоты</w:t></w:r>" +
    "</w:p>"

$screenshotsRange = Get-ParagraphRange("Скриншоты")
$screenshotsRange.InsertXML($screenshotsXml)

Write-Host "Edit applied."
